# Localize handlebars: expose the data model on {{name}} -> {{data.name}}
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("D2").Value = $ws.Range("D2").Value2.Replace("{{name}}", "{{data.name}}")
$ws.Range("D3").Value = $ws.Range("D3").Value2.Replace("{{name}}", "{{data.name}}")
$ws.Range("D4").Value = $ws.Range("D4").Value2.Replace("{{name}}", "{{data.name}}")

# The longer replacement text wraps onto an extra line in row 3
$ws.Rows.Item(3).RowHeight = 60

# Update the selection on the survey sheet
[void]$ws.Range("D10").Select()

# Make "initial" the active/selected sheet (was "settings")
$wsInitial = $wb.Worksheets.Item("initial")
[void]$wsInitial.Select()
